$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 4326
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 4326
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 4326
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -4978

# Row 43
$ws.Range("H43").Value = 2282.923
$ws.Range("I43").Value = 2985.5715
$ws.Range("J43").Value = 1463.1666
$ws.Range("K43").Value = 2985.5715
$ws.Range("L43").Value = 1463.1666
$ws.Range("M43").Value = -2916.5715
$ws.Range("N43").Value = -1601.1666

# Row 88
$ws.Range("H88").Value = 936707.8
$ws.Range("I88").Value = 2852.1667
$ws.Range("J88").Value = 1446083.6
$ws.Range("K88").Value = 2852.1667
$ws.Range("L88").Value = 1446083.6
$ws.Range("M88").Value = -2446.1667
$ws.Range("N88").Value = -1446895.6

# Row 91
$ws.Range("H91").Value = 936707.8
$ws.Range("I91").Value = 2852.1667
$ws.Range("J91").Value = 1446083.6
$ws.Range("K91").Value = 2852.1667
$ws.Range("L91").Value = 1446083.6
$ws.Range("M91").Value = -1448.1667
$ws.Range("N91").Value = -1448891.6

# Row 138
$ws.Range("H138").Value = 1842.39
$ws.Range("I138").Value = 1227.85
$ws.Range("J138").Value = 1996.025
$ws.Range("K138").Value = 3683.55
$ws.Range("L138").Value = 5988.075000000001
$ws.Range("M138").Value = 1456.45
$ws.Range("N138").Value = -16268.075

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 64946.062
$ws.Range("I2").Value = 145496.14
$ws.Range("J2").Value = 2296
$ws.Range("K2").Value = 145496.14
$ws.Range("L2").Value = 2296
$ws.Range("M2").Value = -145383.14
$ws.Range("N2").Value = -2522

# Row 61
$ws.Range("H61").Value = 2166.7778
$ws.Range("I61").Value = 1937.625
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1937.625
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1725.625
$ws.Range("N61").Value = -4424

# Row 88
$ws.Range("H88").Value = 2862.818
$ws.Range("I88").Value = 2399
$ws.Range("J88").Value = 3419.4
$ws.Range("K88").Value = 2399
$ws.Range("L88").Value = 3419.4
$ws.Range("M88").Value = -1993
$ws.Range("N88").Value = -4231.4

# Row 91
$ws.Range("H91").Value = 2862.818
$ws.Range("I91").Value = 2399
$ws.Range("J91").Value = 3419.4
$ws.Range("K91").Value = 2399
$ws.Range("L91").Value = 3419.4
$ws.Range("M91").Value = -995
$ws.Range("N91").Value = -6227.4

# Row 116
$ws.Range("H116").Value = 64946.062
$ws.Range("I116").Value = 145496.14
$ws.Range("J116").Value = 2296
$ws.Range("K116").Value = 145496.14
$ws.Range("L116").Value = 2296
$ws.Range("M116").Value = -143202.14
$ws.Range("N116").Value = -6884

# Row 136
$ws.Range("H136").Value = 2166.7778
$ws.Range("I136").Value = 1937.625
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 5812.875
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -3262.875
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 64946.062
$ws.Range("I3").Value = 145496.14
$ws.Range("J3").Value = 2296
$ws.Range("K3").Value = 145496.14
$ws.Range("L3").Value = 2296
$ws.Range("M3").Value = -145382.14
$ws.Range("N3").Value = -2524

# Row 31
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = ""

# Row 81
$ws.Range("H81").Value = 17296
$ws.Range("J81").Value = 17296
$ws.Range("L81").Value = 17296
$ws.Range("N81").Value = -19418

# Row 84
$ws.Range("H84").Value = 17296
$ws.Range("J84").Value = 17296
$ws.Range("L84").Value = 51888
$ws.Range("N84").Value = -62496

# Row 86
$ws.Range("H86").Value = 2015.0571
$ws.Range("I86").Value = 1766.4688
$ws.Range("J86").Value = 4666.6665
$ws.Range("K86").Value = 1766.4688
$ws.Range("L86").Value = 4666.6665
$ws.Range("M86").Value = -643.4688000000001
$ws.Range("N86").Value = -6912.6665

# Row 89
$ws.Range("H89").Value = 2015.0571
$ws.Range("I89").Value = 1766.4688
$ws.Range("J89").Value = 4666.6665
$ws.Range("K89").Value = 8832.344000000001
$ws.Range("L89").Value = 23333.3325
$ws.Range("M89").Value = -3216.344000000001
$ws.Range("N89").Value = -34565.3325

$ws = $wb.Worksheets.Item("CRP")
# Row 32
$ws.Range("H32").Value = 3339973.2
$ws.Range("I32").Value = 3339973.2
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3339973.2
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3339657.2
$ws.Range("N32").Value = ""

# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = ""

# Row 35
$ws.Range("H35").Value = 758.3333
$ws.Range("I35").Value = 758.3333
$ws.Range("K35").Value = 758.3333
$ws.Range("M35").Value = -464.3333

# Row 36
$ws.Range("H36").Value = 3675
$ws.Range("I36").Value = 950
$ws.Range("J36").Value = 11850
$ws.Range("K36").Value = 950
$ws.Range("L36").Value = 11850
$ws.Range("M36").Value = -562
$ws.Range("N36").Value = -12626

# Row 39
$ws.Range("H39").Value = 15995
$ws.Range("I39").Value = 1990
$ws.Range("J39").Value = 30000
$ws.Range("K39").Value = 1990
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = -1599
$ws.Range("N39").Value = -30782

# Row 40
$ws.Range("H40").Value = 3675
$ws.Range("I40").Value = 950
$ws.Range("J40").Value = 11850
$ws.Range("K40").Value = 950
$ws.Range("L40").Value = 11850
$ws.Range("M40").Value = -790
$ws.Range("N40").Value = -12170

# Row 42
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = ""

# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = ""

# Row 45
$ws.Range("H45").Value = 36689
$ws.Range("I45").Value = 40067
$ws.Range("J45").Value = 35000
$ws.Range("K45").Value = 40067
$ws.Range("L45").Value = 35000
$ws.Range("M45").Value = -39474
$ws.Range("N45").Value = -36186

# Row 49
$ws.Range("H49").Value = 15995
$ws.Range("I49").Value = 1990
$ws.Range("J49").Value = 30000
$ws.Range("K49").Value = 1990
$ws.Range("L49").Value = 30000
$ws.Range("M49").Value = -1808
$ws.Range("N49").Value = -30364

# Row 105
$ws.Range("H105").Value = 2680
$ws.Range("I105").Value = 2440
$ws.Range("K105").Value = 2440
$ws.Range("M105").Value = -693

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1573.909
$ws.Range("J68").Value = 1764.125
$ws.Range("L68").Value = 5292.375
$ws.Range("N68").Value = -6914.375

# Row 71
$ws.Range("H71").Value = 1573.909
$ws.Range("J71").Value = 1764.125
$ws.Range("L71").Value = 15877.125
$ws.Range("N71").Value = -23989.125

# Row 131
$ws.Range("H131").Value = 8790729
$ws.Range("I131").Value = 71573190
$ws.Range("J131").Value = 1183.98
$ws.Range("K131").Value = 214719570
$ws.Range("L131").Value = 3551.94
$ws.Range("M131").Value = -214714530
$ws.Range("N131").Value = -13631.94

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 533.06665
$ws.Range("I107").Value = 359.4
$ws.Range("J107").Value = 619.9
$ws.Range("K107").Value = 359.4
$ws.Range("L107").Value = 619.9
$ws.Range("M107").Value = 1560.6
$ws.Range("N107").Value = -4459.9

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1546.2354
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 1617.875
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 1617.875
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -2207.875

# Row 27
$ws.Range("H27").Value = 1546.2354
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 1617.875
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 1617.875
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -1831.875

# Row 68
$ws.Range("H68").Value = 1935.44
$ws.Range("I68").Value = 2137.7334
$ws.Range("J68").Value = 1632
$ws.Range("K68").Value = 2137.7334
$ws.Range("L68").Value = 1632
$ws.Range("M68").Value = -1388.7334
$ws.Range("N68").Value = -3130

# Row 71
$ws.Range("H71").Value = 1935.44
$ws.Range("I71").Value = 2137.7334
$ws.Range("J71").Value = 1632
$ws.Range("K71").Value = 10688.667
$ws.Range("L71").Value = 8160
$ws.Range("M71").Value = -6944.667000000001
$ws.Range("N71").Value = -15648

# Row 82
$ws.Range("H82").Value = 1693.762
$ws.Range("I82").Value = 1876.7142
$ws.Range("J82").Value = 1327.8572
$ws.Range("K82").Value = 1876.7142
$ws.Range("L82").Value = 1327.8572
$ws.Range("M82").Value = -1515.7142
$ws.Range("N82").Value = -2049.8572

# Row 85
$ws.Range("H85").Value = 1693.762
$ws.Range("I85").Value = 1876.7142
$ws.Range("J85").Value = 1327.8572
$ws.Range("K85").Value = 1876.7142
$ws.Range("L85").Value = 1327.8572
$ws.Range("M85").Value = -628.7141999999999
$ws.Range("N85").Value = -3823.8572

# Row 106
$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524

# Row 132
$ws.Range("H132").Value = 2732.6572
$ws.Range("I132").Value = 2302.2727
$ws.Range("J132").Value = 3461
$ws.Range("K132").Value = 6906.8181
$ws.Range("L132").Value = 10383
$ws.Range("M132").Value = -4376.8181
$ws.Range("N132").Value = -15443

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3530.1
$ws.Range("I62").Value = 4116.8335
$ws.Range("J62").Value = 2650
$ws.Range("K62").Value = 4116.8335
$ws.Range("L62").Value = 2650
$ws.Range("M62").Value = -3492.8335
$ws.Range("N62").Value = -3898

# Row 65
$ws.Range("H65").Value = 3530.1
$ws.Range("I65").Value = 4116.8335
$ws.Range("J65").Value = 2650
$ws.Range("K65").Value = 20584.1675
$ws.Range("L65").Value = 13250
$ws.Range("M65").Value = -17464.1675
$ws.Range("N65").Value = -19490

# Row 132
$ws.Range("H132").Value = 867.17645
$ws.Range("I132").Value = 786
$ws.Range("J132").Value = 1092.6666
$ws.Range("K132").Value = 2358
$ws.Range("L132").Value = 3277.9998
$ws.Range("M132").Value = 172
$ws.Range("N132").Value = -8337.9998
